# Property disclosure workbook: flesh out the "汽車" (car) sheet so its
# rows carry the same name/capacity/owner/.../index column layout used by
# the other property-category sheets (land, building, cash, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Header row -------------------------------------------------------
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Row 2 (car #50) ---------------------------------------------------
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(2, 10).Value = "2011-11-17"
$ws.Cells.Item(2, 11).Value = "吳育昇"
$ws.Cells.Item(2, 12).Value = 1322
$ws.Cells.Item(2, 13).Value = "tmpe6fb1"
$ws.Cells.Item(2, 14).Value = 50

# --- Row 3 (car #51) ---------------------------------------------------
$ws.Cells.Item(3, 8).Value = "land"
$ws.Cells.Item(3, 9).Value = "normal"
$ws.Cells.Item(3, 10).Value = "2011-11-17"
$ws.Cells.Item(3, 11).Value = "吳育昇"
$ws.Cells.Item(3, 12).Value = 1322
$ws.Cells.Item(3, 13).Value = "tmpe6fb1"
$ws.Cells.Item(3, 14).Value = 51

$ws.Range("A1:N3").EntireColumn.AutoFit() | Out-Null
